$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Loc"
$ws.Range("B1").Value = "khongbiet"
$ws.Range("C1").Value = "Loc.png"
$ws.Range("D1").Value = "khongbiet"
$ws.Range("E1").Value = 20
